# Updates the "cryptos" price/volume table with refreshed figures.
# Price cells in column D that are plain decimal numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the workbook's
# existing inline-string cells, e.g. "27.958.45" would otherwise be
# misread as a number) and the cell style is then reset to "Normal" so no
# stray quote-prefix / text-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.958.45'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.859.39'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''311.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '''0.5137'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.89%  '
$ws.Range('D8').Value = '''0.3805'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').Value = '''0.08262'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -9.16%  '
$ws.Range('D10').Value = '''41.61'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').Value = '''6.181'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.859.92'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = '''20.42'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').Value = '''7.194'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '''0.00001091'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = '''90.30'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').Value = '''0.06602'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = '''17.72'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '''5.995'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').Value = '27.989.34'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = '''11.00'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').Value = '''2.217'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.33%  '
$ws.Range('D26').Value = '''2.582'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').Value = '2.074.46'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '''156.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('D30').Value = '''124.35'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('D31').Value = '''0.1062'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = '''1.034'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').Value = '''9.541'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02423'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.06525'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').Value = '''0.2179'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('D39').Value = '''1.205'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '''0.6411'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('D41').Value = '''1.233'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('D42').Value = '''11.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '''4.868'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').Value = '''0.6097'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').Value = '''13.04'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.73%  '
$ws.Range('D46').Value = '''1.281'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '''3.648'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D49').Value = '''1.204'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').Value = '''120.53'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = '''79.46'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.99%  '
